$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1): update 想去人数 (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6770
$ws1.Range("F7").Value = 12
$ws1.Range("F10").Value = 6347
$ws1.Range("F22").Value = 4706
$ws1.Range("F25").Value = 165
$ws1.Range("F26").Value = 197
$ws1.Range("F27").Value = 100

# Sheet "全部类型" (sheetId=4): update 想去人数 (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 6770
$ws4.Range("F7").Value = 12
$ws4.Range("F10").Value = 6347
$ws4.Range("F22").Value = 4706
$ws4.Range("F26").Value = 165
$ws4.Range("F27").Value = 197
$ws4.Range("F28").Value = 100
